$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 3 (source date 44223)
$ws.Range("D2").Value = 44223
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 3500
$ws.Range("L2").Value = 3800
$ws.Range("M2").Value = 3688
$ws.Range("N2").Value = '$/paquete 2 kilos'
$ws.Range("O2").Value = 'Provincia de Diguillín'
$ws.Range("P2").Value = 1844
$ws.Range("Q2").Value = 2

# Row 3 <- original row 7 (source date 44160)
$ws.Range("D3").Value = 44160
$ws.Range("J3").Value = 43
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = 3709
$ws.Range("N3").Value = '$/paquete 36 unidades'
$ws.Range("O3").Value = 'Región Metropolitana'
$ws.Range("P3").Value = 103
$ws.Range("Q3").Value = 36

# Row 4 <- original row 10 (source date 44215)
$ws.Range("D4").Value = 44215
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 3500
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = 3768
$ws.Range("N4").Value = '$/paquete 2 kilos'
$ws.Range("O4").Value = 'Provincia de Diguillín'
$ws.Range("P4").Value = 1884
$ws.Range("Q4").Value = 2

# Row 5 <- original row 8 (source date 44208)
$ws.Range("D5").Value = 44208
$ws.Range("J5").Value = 85
$ws.Range("K5").Value = 3700
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 3824
$ws.Range("N5").Value = '$/paquete 2 kilos'
$ws.Range("O5").Value = 'Provincia de Diguillín'
$ws.Range("P5").Value = 1912
$ws.Range("Q5").Value = 2

# Row 6 <- original row 4 (source date 44209)
$ws.Range("D6").Value = 44209
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 3500
$ws.Range("L6").Value = 4000
$ws.Range("M6").Value = 3767
$ws.Range("N6").Value = '$/paquete 2 kilos'
$ws.Range("O6").Value = 'Provincia de Diguillín'
$ws.Range("P6").Value = 1884
$ws.Range("Q6").Value = 2

# Row 7 <- original row 5 (source date 44210)
$ws.Range("D7").Value = 44210
$ws.Range("J7").Value = 105
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = 3714
$ws.Range("N7").Value = '$/paquete 2 kilos'
$ws.Range("O7").Value = 'Provincia de Diguillín'
$ws.Range("P7").Value = 1857
$ws.Range("Q7").Value = 2

# Row 8 <- original row 9 (source date 44166)
$ws.Range("D8").Value = 44166
$ws.Range("J8").Value = 70
$ws.Range("K8").Value = 3500
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 3679
$ws.Range("N8").Value = '$/paquete 36 unidades'
$ws.Range("O8").Value = 'Región Metropolitana'
$ws.Range("P8").Value = 102
$ws.Range("Q8").Value = 36

# Row 9 <- original row 2 (source date 44225)
$ws.Range("D9").Value = 44225
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 3400
$ws.Range("L9").Value = 3700
$ws.Range("M9").Value = 3550
$ws.Range("N9").Value = '$/paquete 2 kilos'
$ws.Range("O9").Value = 'Provincia de Diguillín'
$ws.Range("P9").Value = 1775
$ws.Range("Q9").Value = 2

# Row 10 <- original row 6 (source date 44161)
$ws.Range("D10").Value = 44161
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 2800
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2900
$ws.Range("N10").Value = '$/paquete 2 kilos'
$ws.Range("O10").Value = 'Provincia de Diguillín'
$ws.Range("P10").Value = 1450
$ws.Range("Q10").Value = 2

